$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(45, 8).Value = 3000  # H45: was 1966.3334
$ws.Cells.Item(45, 9).Value = 0  # I45: was 1966.3334
$ws.Cells.Item(45, 10).Value = 3000  # J45: was 0
$ws.Cells.Item(45, 11).Value = 0  # K45: was 5899.0002
$ws.Cells.Item(45, 12).Value = 9000  # L45: was 0
$ws.Cells.Item(45, 13).ClearContents()  # M45: was -5707.0002
$ws.Cells.Item(45, 14).Value = -9384  # N45: was None
$ws.Cells.Item(124, 8).Value = 32597.5  # H124: was 29950
$ws.Cells.Item(124, 10).Value = 32597.5  # J124: was 29950
$ws.Cells.Item(124, 12).Value = 32597.5  # L124: was 29950
$ws.Cells.Item(124, 14).Value = -42417.5  # N124: was -39770

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 2983.3333  # H4: was 2683.6667
$ws.Cells.Item(4, 9).Value = 0  # I4: was 101
$ws.Cells.Item(4, 10).Value = 2983.3333  # J4: was 3975
$ws.Cells.Item(4, 11).Value = 0  # K4: was 101
$ws.Cells.Item(4, 12).Value = 2983.3333  # L4: was 3975
$ws.Cells.Item(4, 13).ClearContents()  # M4: was 15
$ws.Cells.Item(4, 14).Value = -3215.3333  # N4: was -4207
$ws.Cells.Item(5, 8).Value = 137.5  # H5: was 150
$ws.Cells.Item(5, 10).Value = 0  # J5: was 200
$ws.Cells.Item(5, 12).Value = 0  # L5: was 200
$ws.Cells.Item(5, 14).ClearContents()  # N5: was -424
$ws.Cells.Item(6, 8).Value = 3000  # H6: was 1000000
$ws.Cells.Item(6, 9).Value = 3000  # I6: was 0
$ws.Cells.Item(6, 10).Value = 0  # J6: was 1000000
$ws.Cells.Item(6, 11).Value = 3000  # K6: was 0
$ws.Cells.Item(6, 12).Value = 0  # L6: was 1000000
$ws.Cells.Item(6, 13).Value = -2827  # M6: was None
$ws.Cells.Item(6, 14).ClearContents()  # N6: was -1000346
$ws.Cells.Item(8, 8).Value = 0  # H8: was 5
$ws.Cells.Item(8, 9).Value = 0  # I8: was 5
$ws.Cells.Item(8, 11).Value = 0  # K8: was 5
$ws.Cells.Item(8, 13).ClearContents()  # M8: was 139
$ws.Cells.Item(16, 8).Value = 0  # H16: was 500
$ws.Cells.Item(16, 9).Value = 0  # I16: was 500
$ws.Cells.Item(16, 11).Value = 0  # K16: was 500
$ws.Cells.Item(16, 13).ClearContents()  # M16: was -213
$ws.Cells.Item(74, 8).Value = 6098169.5  # H74: was 6449968
$ws.Cells.Item(74, 9).Value = 7961182.5  # I74: was 8359225.5
$ws.Cells.Item(74, 10).Value = 79203.92  # J74: was 85775.914
$ws.Cells.Item(74, 11).Value = 7961182.5  # K74: was 8359225.5
$ws.Cells.Item(74, 12).Value = 79203.92  # L74: was 85775.914
$ws.Cells.Item(74, 13).Value = -7960308.5  # M74: was -8358351.5
$ws.Cells.Item(74, 14).Value = -80951.92  # N74: was -87523.914
$ws.Cells.Item(77, 8).Value = 6098169.5  # H77: was 6449968
$ws.Cells.Item(77, 9).Value = 7961182.5  # I77: was 8359225.5
$ws.Cells.Item(77, 10).Value = 79203.92  # J77: was 85775.914
$ws.Cells.Item(77, 11).Value = 39805912.5  # K77: was 41796127.5
$ws.Cells.Item(77, 12).Value = 396019.6  # L77: was 428879.57
$ws.Cells.Item(77, 13).Value = -39801544.5  # M77: was -41791759.5
$ws.Cells.Item(77, 14).Value = -404755.6  # N77: was -437615.57
$ws.Cells.Item(132, 8).Value = 39148.98  # H132: was 38959.035
$ws.Cells.Item(132, 9).Value = 23768.777  # I132: was 24136.477
$ws.Cells.Item(132, 10).Value = 102068  # J132: was 93308.414
$ws.Cells.Item(132, 11).Value = 71306.33099999999  # K132: was 72409.431
$ws.Cells.Item(132, 12).Value = 306204  # L132: was 279925.242
$ws.Cells.Item(132, 13).Value = -68776.33099999999  # M132: was -69879.431
$ws.Cells.Item(132, 14).Value = -311264  # N132: was -284985.242

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 137.5  # H4: was 150
$ws.Cells.Item(4, 10).Value = 0  # J4: was 200
$ws.Cells.Item(4, 12).Value = 0  # L4: was 200
$ws.Cells.Item(4, 14).ClearContents()  # N4: was -430

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6275.8066  # H31: was 6511.5
$ws.Cells.Item(31, 9).Value = 2366.6667  # I31: was 2517.7273
$ws.Cells.Item(31, 10).Value = 8744.736999999999  # J31: was 8823.684999999999
$ws.Cells.Item(31, 11).Value = 2366.6667  # K31: was 2517.7273
$ws.Cells.Item(31, 12).Value = 8744.736999999999  # L31: was 8823.684999999999
$ws.Cells.Item(31, 13).Value = -2071.6667  # M31: was -2222.7273
$ws.Cells.Item(31, 14).Value = -9334.736999999999  # N31: was -9413.684999999999
$ws.Cells.Item(34, 8).Value = 6275.8066  # H34: was 6511.5
$ws.Cells.Item(34, 9).Value = 2366.6667  # I34: was 2517.7273
$ws.Cells.Item(34, 10).Value = 8744.736999999999  # J34: was 8823.684999999999
$ws.Cells.Item(34, 11).Value = 2366.6667  # K34: was 2517.7273
$ws.Cells.Item(34, 12).Value = 8744.736999999999  # L34: was 8823.684999999999
$ws.Cells.Item(34, 13).Value = -2164.6667  # M34: was -2315.7273
$ws.Cells.Item(34, 14).Value = -9148.736999999999  # N34: was -9227.684999999999
$ws.Cells.Item(119, 8).Value = 48238.75  # H119: was 0
$ws.Cells.Item(119, 10).Value = 48238.75  # J119: was 0
$ws.Cells.Item(119, 12).Value = 48238.75  # L119: was 0
$ws.Cells.Item(119, 14).Value = -57914.75  # N119: was None

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(49, 8).Value = 2500  # H49: was 1862.5
$ws.Cells.Item(49, 9).Value = 0  # I49: was 400
$ws.Cells.Item(49, 10).Value = 2500  # J49: was 2071.4285
$ws.Cells.Item(49, 11).Value = 0  # K49: was 1200
$ws.Cells.Item(49, 12).Value = 7500  # L49: was 6214.2855
$ws.Cells.Item(49, 13).ClearContents()  # M49: was -1044
$ws.Cells.Item(49, 14).Value = -7812  # N49: was -6526.2855
$ws.Cells.Item(129, 8).Value = 3790127.5  # H129: was 3625374.2
$ws.Cells.Item(129, 9).Value = 1913.5  # I129: was 1789.7778
$ws.Cells.Item(129, 11).Value = 5740.5  # K129: was 5369.3334
$ws.Cells.Item(129, 13).Value = -740.5  # M129: was -369.3334000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 69.64286  # H2: was 65.8
$ws.Cells.Item(2, 10).Value = 122.28571  # J2: was 108.5
$ws.Cells.Item(2, 12).Value = 122.28571  # L2: was 108.5
$ws.Cells.Item(2, 14).Value = -348.28571  # N2: was -334.5
$ws.Cells.Item(48, 8).Value = 0  # H48: was 8000
$ws.Cells.Item(48, 9).Value = 0  # I48: was 8000
$ws.Cells.Item(48, 11).Value = 0  # K48: was 8000
$ws.Cells.Item(48, 13).ClearContents()  # M48: was -7515
$ws.Cells.Item(80, 8).Value = 3258.1072  # H80: was 3324.1428
$ws.Cells.Item(80, 9).Value = 2782.2727  # I80: was 2823.0908
$ws.Cells.Item(80, 10).Value = 3566  # J80: was 3648.353
$ws.Cells.Item(80, 11).Value = 2782.2727  # K80: was 2823.0908
$ws.Cells.Item(80, 12).Value = 3566  # L80: was 3648.353
$ws.Cells.Item(80, 13).Value = -1784.2727  # M80: was -1825.0908
$ws.Cells.Item(80, 14).Value = -5562  # N80: was -5644.353
$ws.Cells.Item(83, 8).Value = 3258.1072  # H83: was 3324.1428
$ws.Cells.Item(83, 9).Value = 2782.2727  # I83: was 2823.0908
$ws.Cells.Item(83, 10).Value = 3566  # J83: was 3648.353
$ws.Cells.Item(83, 11).Value = 13911.3635  # K83: was 14115.454
$ws.Cells.Item(83, 12).Value = 17830  # L83: was 18241.765
$ws.Cells.Item(83, 13).Value = -8919.363499999999  # M83: was -9123.454
$ws.Cells.Item(83, 14).Value = -27814  # N83: was -28225.765

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(47, 8).Value = 9994  # H47: was 10960
$ws.Cells.Item(47, 10).Value = 9994  # J47: was 10960
$ws.Cells.Item(47, 12).Value = 9994  # L47: was 10960
$ws.Cells.Item(47, 14).Value = -10974  # N47: was -11940
$ws.Cells.Item(52, 8).Value = 9994  # H52: was 10960
$ws.Cells.Item(52, 10).Value = 9994  # J52: was 10960
$ws.Cells.Item(52, 12).Value = 9994  # L52: was 10960
$ws.Cells.Item(52, 14).Value = -10460  # N52: was -11426
$ws.Cells.Item(82, 8).Value = 1120.1818  # H82: was 1296.8889
$ws.Cells.Item(82, 9).Value = 912.75  # I82: was 1059.6666
$ws.Cells.Item(82, 10).Value = 1238.7142  # J82: was 1415.5
$ws.Cells.Item(82, 11).Value = 912.75  # K82: was 1059.6666
$ws.Cells.Item(82, 12).Value = 1238.7142  # L82: was 1415.5
$ws.Cells.Item(82, 13).Value = -551.75  # M82: was -698.6666
$ws.Cells.Item(82, 14).Value = -1960.7142  # N82: was -2137.5
$ws.Cells.Item(85, 8).Value = 1120.1818  # H85: was 1296.8889
$ws.Cells.Item(85, 9).Value = 912.75  # I85: was 1059.6666
$ws.Cells.Item(85, 10).Value = 1238.7142  # J85: was 1415.5
$ws.Cells.Item(85, 11).Value = 912.75  # K85: was 1059.6666
$ws.Cells.Item(85, 12).Value = 1238.7142  # L85: was 1415.5
$ws.Cells.Item(85, 13).Value = 335.25  # M85: was 188.3334
$ws.Cells.Item(85, 14).Value = -3734.7142  # N85: was -3911.5
$ws.Cells.Item(93, 8).Value = 1099.2858  # H93: was 1262.1111
$ws.Cells.Item(93, 9).Value = 1099.2858  # I93: was 1107.375
$ws.Cells.Item(93, 10).Value = 0  # J93: was 2500
$ws.Cells.Item(93, 11).Value = 1099.2858  # K93: was 1107.375
$ws.Cells.Item(93, 12).Value = 0  # L93: was 2500
$ws.Cells.Item(93, 13).Value = 148.7141999999999  # M93: was 140.625
$ws.Cells.Item(93, 14).ClearContents()  # N93: was -4996
$ws.Cells.Item(100, 8).Value = 1710.9375  # H100: was 1799.0667
$ws.Cells.Item(100, 9).Value = 1538.5  # I100: was 1625.75
$ws.Cells.Item(100, 10).Value = 1998.3334  # J100: was 1997.1428
$ws.Cells.Item(100, 11).Value = 1538.5  # K100: was 1625.75
$ws.Cells.Item(100, 12).Value = 1998.3334  # L100: was 1997.1428
$ws.Cells.Item(100, 13).Value = -997.5  # M100: was -1084.75
$ws.Cells.Item(100, 14).Value = -3080.3334  # N100: was -3079.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 0  # H47: was 8750
$ws.Cells.Item(47, 10).Value = 0  # J47: was 8750
$ws.Cells.Item(47, 12).Value = 0  # L47: was 8750
$ws.Cells.Item(47, 14).ClearContents()  # N47: was -9894
$ws.Cells.Item(81, 8).Value = 1847.15  # H81: was 2041.1052
$ws.Cells.Item(81, 9).Value = 737.7778  # I81: was 754
$ws.Cells.Item(81, 10).Value = 2754.818  # J81: was 2791.9167
$ws.Cells.Item(81, 11).Value = 1475.5556  # K81: was 1508
$ws.Cells.Item(81, 12).Value = 5509.636  # L81: was 5583.8334
$ws.Cells.Item(81, 13).Value = -414.5555999999999  # M81: was -447
$ws.Cells.Item(81, 14).Value = -7631.636  # N81: was -7705.8334
$ws.Cells.Item(84, 8).Value = 1847.15  # H84: was 2041.1052
$ws.Cells.Item(84, 9).Value = 737.7778  # I84: was 754
$ws.Cells.Item(84, 10).Value = 2754.818  # J84: was 2791.9167
$ws.Cells.Item(84, 11).Value = 7377.777999999999  # K84: was 7540
$ws.Cells.Item(84, 12).Value = 27548.18  # L84: was 27919.167
$ws.Cells.Item(84, 13).Value = -2073.777999999999  # M84: was -2236
$ws.Cells.Item(84, 14).Value = -38156.18  # N84: was -38527.167
$ws.Cells.Item(121, 8).Value = 60420  # H121: was 0
$ws.Cells.Item(121, 10).Value = 60420  # J121: was 0
$ws.Cells.Item(121, 12).Value = 60420  # L121: was 0
$ws.Cells.Item(121, 14).Value = -63914  # N121: was None
